$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 44358
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("P2").Value = 1450
$ws.Range("D4").Value = 44265
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 1550
$ws.Range("D5").Value = 44441
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15500
$ws.Range("P5").Value = 1550
$ws.Range("D6").Value = 44428
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 1550
$ws.Range("D7").Value = 44918
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 12250
$ws.Range("P7").Value = 1225
$ws.Range("D8").Value = 44263
$ws.Range("D9").Value = 44943
$ws.Range("D10").Value = 44547
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19500
$ws.Range("P10").Value = 1950
$ws.Range("D11").Value = 44330
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("P11").Value = 1350
$ws.Range("D12").Value = 44727
$ws.Range("J12").Value = 400
$ws.Range("L12").Value = 19000
$ws.Range("M12").Value = 18500
$ws.Range("P12").Value = 1850
$ws.Range("D13").Value = 44847
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 16500
$ws.Range("P13").Value = 1650
$ws.Range("D14").Value = 44893
$ws.Range("J14").Value = 1400
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15571
$ws.Range("P14").Value = 1557
$ws.Range("D15").Value = 44942
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("P15").Value = 1450
$ws.Range("D16").Value = 44406
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("P16").Value = 1450
$ws.Range("D17").Value = 44860
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("P17").Value = 1450
$ws.Range("D18").Value = 44890
$ws.Range("K18").Value = 16000
$ws.Range("L18").Value = 17000
$ws.Range("M18").Value = 16500
$ws.Range("P18").Value = 1650
$ws.Range("D19").Value = 44460
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15500
$ws.Range("P19").Value = 1550
$ws.Range("D20").Value = 44679
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 19000
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = 19500
$ws.Range("P20").Value = 1950
$ws.Range("D21").Value = 44218
$ws.Range("J21").Value = 320
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 11000
$ws.Range("M21").Value = 10500
$ws.Range("P21").Value = 1050
$ws.Range("D22").Value = 44580
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 19000
$ws.Range("P22").Value = 1900
$ws.Range("D23").Value = 44714
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 19000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 19500
$ws.Range("P23").Value = 1950
$ws.Range("D24").Value = 44377
$ws.Range("J24").Value = 650
$ws.Range("M24").Value = 14538
$ws.Range("P24").Value = 1454
$ws.Range("D25").Value = 44160
$ws.Range("J25").Value = 360
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 11000
$ws.Range("M25").Value = 10500
$ws.Range("P25").Value = 1050
$ws.Range("D27").Value = 44914
$ws.Range("J27").Value = 100
$ws.Range("D28").Value = 44204
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 11000
$ws.Range("M28").Value = 10500
$ws.Range("P28").Value = 1050
$ws.Range("D29").Value = 44524
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 21000
$ws.Range("M29").Value = 20500
$ws.Range("P29").Value = 2050
$ws.Range("D30").Value = 44882
$ws.Range("J30").Value = 400
$ws.Range("M30").Value = 15550
$ws.Range("P30").Value = 1555
$ws.Range("D31").Value = 44694
$ws.Range("J31").Value = 400
$ws.Range("K31").Value = 16000
$ws.Range("L31").Value = 17000
$ws.Range("M31").Value = 16500
$ws.Range("P31").Value = 1650
$ws.Range("D32").Value = 44291
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13500
$ws.Range("P32").Value = 1350
$ws.Range("D33").Value = 44904
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 14000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 14500
$ws.Range("P33").Value = 1450
$ws.Range("D34").Value = 44644
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 20000
$ws.Range("L34").Value = 21000
$ws.Range("M34").Value = 20500
$ws.Range("P34").Value = 2050
